# 6.4.2.1 Общий объем забора пресной воды
# Adds a new "2023" data column (M) to the existing yearly table (2014-2022)
# and gives row 4 (the year header row) localized column headers in
# A4/B4/C4 that mirror the ones already used in A1:C2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 header labels (were blank) ------------------------------------
$ws.Range("A4").Value2 = "Көрсөткүчтөрдүн аталыштары"
$ws.Range("B4").Value2 = "Наименование показателей"
$ws.Range("C4").Value2 = "Items"

# --- New column M (year 2023) data ----------------------------------------
$ws.Range("M4").Value2  = 2023
$ws.Range("M5").Value2  = 8872.5
$ws.Range("M7").Value2  = 8601.5
$ws.Range("M8").Value2  = 271
$ws.Range("M10").Value2 = 723.4
$ws.Range("M11").Value2 = 1205.5999999999999
$ws.Range("M12").Value2 = 779.6
$ws.Range("M13").Value2 = 829.3
$ws.Range("M14").Value2 = 1314.9
$ws.Range("M15").Value2 = 1034.5999999999999
$ws.Range("M16").Value2 = 2762.1
$ws.Range("M17").Value2 = 166
$ws.Range("M18").Value2 = 57

# M6 / M9 are blank "section" rows -- mirror their row's formatting only,
# no value to set.

# --- Carry over formatting from column L (same row) into column M --------
$rows = 4..18
foreach ($r in $rows) {
    $src = $ws.Cells.Item($r, 12)  # L
    $dst = $ws.Cells.Item($r, 13)  # M
    $dst.Font.Name = $src.Font.Name
    $dst.Font.Bold = $src.Font.Bold
    $dst.Font.Italic = $src.Font.Italic
    $dst.Font.Size = $src.Font.Size
    $dst.NumberFormat = $src.NumberFormat
    $dst.HorizontalAlignment = $src.HorizontalAlignment
    $dst.VerticalAlignment = $src.VerticalAlignment
    $dst.WrapText = $src.WrapText
}

# Row 4 and row 18 carry a border (medium rule under the header, medium
# rule under the totals) -- replicate it on the new M cells.
$ws.Range("M4").Borders.Item(8).LineStyle = 1
$ws.Range("M4").Borders.Item(8).Weight = -4138
$ws.Range("M4").Borders.Item(9).LineStyle = 1
$ws.Range("M4").Borders.Item(9).Weight = -4138

$ws.Range("M18").Borders.Item(9).LineStyle = 1
$ws.Range("M18").Borders.Item(9).Weight = -4138

# --- Row heights: every data row now gets an explicit 14.25pt height ------
$heights = @{
    4 = 14.25; 5 = 14.25; 6 = 14.25; 7 = 14.25; 8 = 14.25; 9 = 14.25;
    10 = 14.25; 11 = 14.25; 12 = 14.25; 13 = 14.25; 14 = 14.25;
    15 = 14.25; 16 = 14.25; 17 = 14.25; 18 = 14.25
}
foreach ($r in $heights.Keys) {
    $ws.Rows($r).RowHeight = $heights[$r]
}

# --- Drop the stray O2 selection marker on the sheet view -----------------
$ws.Range("A1").Select()
